$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.557.77"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "3.981.55"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "588.31"
$ws.Range("E5").Value = "  +12.23%  "
$ws.Range("D6").Value = "151.97"
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("D7").Value = "0.677"
$ws.Range("E7").Value = "  -2.87%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.746"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").Value = "53.05"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "0.0000317"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "10.76"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "4.618.99"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "3.989.92"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  +8.53%  "
$ws.Range("D17").Value = "13.97"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "20.41"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "72.511.95"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "428.70"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  +13.94%  "
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "3.43"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "4.51"
$ws.Range("E25").Value = "  +22.18%  "
$ws.Range("D26").Value = "14.21"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "11.23"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "10.53"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "36.24"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").Value = "7.81"
$ws.Range("E31").Value = "  +5.71%  "
$ws.Range("D32").Value = "49.99"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D33").Value = "13.44"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "680.71"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "68.44"
$ws.Range("E36").Value = "  +4.47%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "0.0₃0850"
$ws.Range("E38").Value = "  +3.70%  "
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +1.09%  "
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "11.09"
$ws.Range("E42").Value = "  +11.65%  "
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "0.0485"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "3.44"
$ws.Range("E49").Value = "  +5.88%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "2.14"
$ws.Range("E51").Value = "  +6.69%  "
